# Update "Horarios Línea 141" workbook (scrape run at 16:51:51).
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- Sheet "LP1912": header banner + row count ---
$ws1.Range("A2").Value = "Última actualización: 16:51:51"
$ws1.Range("A3").Value = "Total filas: 268"

# Two arrivals recorded at the same scrape/arrival time got their "Linea"
# swapped when the source re-sorted ties - fix the Linea (column C) values.
$ws1.Cells.Item(38, 3).Value = "15_ABASTO"
$ws1.Cells.Item(39, 3).Value = "11_ETCHEVERRY"

$ws1.Cells.Item(106, 3).Value = "10_OLMOS"
$ws1.Cells.Item(107, 3).Value = "16_SANTA ANA"

# Two rows whose ordering (and therefore full data) swapped.
$ws1.Cells.Item(139, 1).Value = "10:49:38"
$ws1.Cells.Item(139, 3).Value = "27_EL RETIRO"
$ws1.Cells.Item(139, 4).Value = 107

$ws1.Cells.Item(140, 1).Value = "11:53:44"
$ws1.Cells.Item(140, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(140, 4).Value = 43

$ws1.Cells.Item(190, 1).Value = "13:55:43"
$ws1.Cells.Item(190, 3).Value = "16_SANTA ANA"
$ws1.Cells.Item(190, 4).Value = 56

$ws1.Cells.Item(191, 1).Value = "13:41:21"
$ws1.Cells.Item(191, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(191, 4).Value = 70

# New arrival scraped at 16:51:51 inserted before the former row 255,
# pushing the rest of the tail (old rows 255-272) down to 256-273.
$ws1.Rows.Item(255).Insert()
$ws1.Cells.Item(255, 1).Value = "16:51:51"
$ws1.Cells.Item(255, 2).Value = "17:53"
$ws1.Cells.Item(255, 3).Value = "10_OLMOS"
$ws1.Cells.Item(255, 4).Value = 62
$ws1.Cells.Item(255, 5).Value = "LP1912"

# --- Sheets "LP1912-215" and "6203-6173": header banner only ---
$ws2.Range("A2").Value = "Última actualización: 16:51:51"
$ws3.Range("A2").Value = "Última actualización: 16:51:51"
